# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the bespoke "Table_0" style ({2A9BDF65-6696-44FB-8456-168A9494A209})
#    to the built-in gallery style {70DEFBB0-8C00-4821-A23B-0D8BEC48BD77}.
#
# 2) The deck's colour theme is recoloured from the "Red Violet"/"Integral"
#    palette to the stock "Office" palette (dk1/lt1 stay black/white; the
#    other ten slots change). Font scheme and format scheme are already
#    identical between the two theme parts, so only the colour slots need
#    to move.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------
$newTableStyle = "{70DEFBB0-8C00-4821-A23B-0D8BEC48BD77}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Recolour the theme back to the stock "Office" palette --------
# MsoThemeColorSchemeIndex slots: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. RGB is written VBA-style (R + G*256 + B*65536).
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000 (unchanged)
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF (unchanged)
$tcs.Item(3).RGB  = 0x6A546B
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5 (unchanged)
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
